$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# --- Update timestamp text in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 21:20"

# --- Swap Asturias / Malaga rows (18 & 19), with Asturias getting new counts ---
# Row 18 becomes Asturias with updated numbers
$ws.Range("A18").Value = "Asturias"
$ws.Range("B18").Value = 1158
$ws.Range("C18").Value = 78
$ws.Range("D18").Value = 1032
$ws.Range("E18").Value = 48

# Row 19 becomes Malaga, carrying the old Malaga numbers
$ws.Range("A19").Value = "Malaga"
$ws.Range("B19").Value = 1158
$ws.Range("C19").Value = 81
$ws.Range("D19").Value = 1013
$ws.Range("E19").Value = 64

# --- Swap Ibiza / Ceuta rows (58 & 59), with Ceuta getting new counts ---
# Row 58 becomes Ceuta with updated numbers
$ws.Range("A58").Value = "Ceuta"
$ws.Range("B58").Value = 27
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 26
$ws.Range("E58").Value = 1

# Row 59 becomes Ibiza, carrying the old Ibiza numbers
$ws.Range("A59").Value = "Ibiza"
$ws.Range("B59").Value = 21
$ws.Range("C59").Value = 18
$ws.Range("D59").Value = 20
$ws.Range("E59").Value = 1
